$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: CageNumber (F10) changes from 2 to 3
$ws.Cells.Item(10, 6).Value = 3

# Row 11: update to new bird values
$ws.Cells.Item(11, 1).Value = 11
$ws.Cells.Item(11, 2).Value = "Golden Australian "
$ws.Cells.Item(11, 3).Value = "Coastal cities"
$ws.Cells.Item(11, 4).Value = 45047
$ws.Cells.Item(11, 5).Value = "Male"
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(11, 7).Value = 3
$ws.Cells.Item(11, 8).Value = 3

# Row 12: update to new bird values
$ws.Cells.Item(12, 1).Value = 12
$ws.Cells.Item(12, 2).Value = "Golden Australian "
$ws.Cells.Item(12, 3).Value = "Central Australia"
$ws.Cells.Item(12, 4).Value = 45075
$ws.Cells.Item(12, 5).Value = "Male"
$ws.Cells.Item(12, 6).Value = 6
$ws.Cells.Item(12, 7).Value = 5
$ws.Cells.Item(12, 8).Value = 4

# Copy date formatting (column D) from an existing formatted cell down to the new rows
$ws.Cells.Item(2, 4).Copy()
$ws.Range("D13:D17").PasteSpecial(-4122)

# Row 13: new bird
$ws.Cells.Item(13, 1).Value = 13
$ws.Cells.Item(13, 2).Value = "Golden European"
$ws.Cells.Item(13, 3).Value = "West Europe"
$ws.Cells.Item(13, 4).Value = 45074
$ws.Cells.Item(13, 5).Value = "Male"
$ws.Cells.Item(13, 6).Value = 3
$ws.Cells.Item(13, 7).Value = 2
$ws.Cells.Item(13, 8).Value = 1

# Row 14: new bird
$ws.Cells.Item(14, 1).Value = 14
$ws.Cells.Item(14, 2).Value = "Golden European"
$ws.Cells.Item(14, 3).Value = "East Europe"
$ws.Cells.Item(14, 4).Value = 45082
$ws.Cells.Item(14, 5).Value = "Female"
$ws.Cells.Item(14, 6).Value = 4
$ws.Cells.Item(14, 7).Value = 12
$ws.Cells.Item(14, 8).Value = 12

# Row 15: new bird (no Gender value in this row)
$ws.Cells.Item(15, 1).Value = 15
$ws.Cells.Item(15, 2).Value = "Golden European"
$ws.Cells.Item(15, 3).Value = "West Europe"
$ws.Cells.Item(15, 4).Value = 45078
$ws.Cells.Item(15, 6).Value = 3
$ws.Cells.Item(15, 7).Value = 13
$ws.Cells.Item(15, 8).Value = 4

# Row 16: new bird
$ws.Cells.Item(16, 1).Value = 16
$ws.Cells.Item(16, 2).Value = "Golden Australian "
$ws.Cells.Item(16, 3).Value = "Coastal cities"
$ws.Cells.Item(16, 4).Value = 45083
$ws.Cells.Item(16, 5).Value = "Male"
$ws.Cells.Item(16, 6).Value = 2
$ws.Cells.Item(16, 7).Value = 11
$ws.Cells.Item(16, 8).Value = 10

# Row 17: new bird
$ws.Cells.Item(17, 1).Value = 17
$ws.Cells.Item(17, 2).Value = "Golden European"
$ws.Cells.Item(17, 3).Value = "West Europe"
$ws.Cells.Item(17, 4).Value = 45077
$ws.Cells.Item(17, 5).Value = "Male"
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 2
$ws.Cells.Item(17, 8).Value = 2

# Refresh the worksheet sort state to cover the newly added rows
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($ws.Range("A2:H17"))
$ws.Sort.Header = 0
$ws.Sort.Apply()
